# Natmi re-run following Dr Hou's advice: the Wnt5a -> Fzd6 LR-pair sheet is
# recomputed with an extra "M2" cluster in the mix (and the sending/target
# cluster combinations revisited), so rows 2-3 change in place and 6 brand
# new rows (4-9) are appended, covering every Sending x Target cluster
# combination for FAPs/sCs as senders across FAPs/sCs/ECs/M2 as targets.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns, in order: A..T
#   Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#   Ligand-expressing cells, Ligand detection rate,
#   Ligand average expression value, Ligand total expression value,
#   Ligand derived specificity of average/total expression value,
#   Receptor-expressing cells, Receptor detection rate,
#   Receptor average expression value, Receptor total expression value,
#   Receptor derived specificity of average/total expression value,
#   Edge average expression weight, Edge total expression weight,
#   Edge average/total expression derived specificity
$data = @(
    ,('FAPs', 'Wnt5a', 'Fzd6', 'ECs', 3, 1, 10.43424333333333, 31.30273, 0.9711091978791583, 0.9711091978791584, 2, 0.6666666666666666, 13.21223933333333, 39.636718, 0.8149747101495924, 0.8149747101495927, 137.8597201822378, 1240.73748164014, 0.7914294370651702, 0.7914294370651705)
    ,('FAPs', 'Wnt5a', 'Fzd6', 'FAPs', 3, 1, 10.43424333333333, 31.30273, 0.9711091978791583, 0.9711091978791584, 3, 1, 2.662736333333334, 7.988209000000001, 0.1642464018940561, 0.1642464018940561, 27.78363883450778, 250.05274951057, 0.1595011915978747, 0.1595011915978747)
    ,('FAPs', 'Wnt5a', 'Fzd6', 'M2', 3, 1, 10.43424333333333, 31.30273, 0.9711091978791583, 0.9711091978791584, 1, 0.3333333333333333, 0.009795, 0.029385, 0.0006041880626379251, 0.0006041880626379253, 0.10220341345, 0.91983072105, 0.0005867325848764781, 0.0005867325848764784)
    ,('FAPs', 'Wnt5a', 'Fzd6', 'sCs', 3, 1, 10.43424333333333, 31.30273, 0.9711091978791583, 0.9711091978791584, 3, 1, 0.327069, 0.9812069999999999, 0.02017469989371348, 0.02017469989371348, 3.41271753279, 30.71445779510999, 0.01959183663123683, 0.01959183663123684)
    ,('sCs', 'Wnt5a', 'Fzd6', 'ECs', 2, 0.6666666666666666, 0.310422, 0.9312659999999999, 0.02889080212084161, 0.02889080212084161, 2, 0.6666666666666666, 13.21223933333333, 39.636718, 0.8149747101495924, 0.8149747101495927, 4.101369758332, 36.912327824988, 0.02354527308442212, 0.02354527308442213)
    ,('sCs', 'Wnt5a', 'Fzd6', 'FAPs', 2, 0.6666666666666666, 0.310422, 0.9312659999999999, 0.02889080212084161, 0.02889080212084161, 3, 1, 2.662736333333334, 7.988209000000001, 0.1642464018940561, 0.1642464018940561, 0.826571938066, 7.439147442594001, 0.004745210296181398, 0.0047452102961814)
    ,('sCs', 'Wnt5a', 'Fzd6', 'M2', 2, 0.6666666666666666, 0.310422, 0.9312659999999999, 0.02889080212084161, 0.02889080212084161, 1, 0.3333333333333333, 0.009795, 0.029385, 0.0006041880626379251, 0.0006041880626379253, 0.00304058349, 0.02736525141, 0.000017455477761446948901154913, 0.00001745547776144695906555028)
    ,('sCs', 'Wnt5a', 'Fzd6', 'sCs', 2, 0.6666666666666666, 0.310422, 0.9312659999999999, 0.02889080212084161, 0.02889080212084161, 3, 1, 0.327069, 0.9812069999999999, 0.02017469989371348, 0.02017469989371348, 0.101529413118, 0.9137647180619999, 0.0005828632624766403, 0.0005828632624766406)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

